# Lines_IPR.xlsx: update the "cant_line" (column D) values for rows 2, 4,
# 5 and 7, and give those four cells a dedicated cell style (same Arial
# font as before, now saved with an explicit-but-default alignment block)
# instead of the plain style they previously shared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-LineCount([string]$addr, [double]$value) {
    $cell = $ws.Range($addr)
    $cell.Value = $value
    # Touching an alignment property (even to its existing default) makes
    # the engine materialise a new, distinct cellXf for this cell -
    # matching the new 5th entry added to cellXfs in the target file.
    $cell.WrapText = $false
}

Set-LineCount "D2" 4
Set-LineCount "D4" 4
Set-LineCount "D5" 3
Set-LineCount "D7" 4

Write-Output "updated D2, D4, D5, D7"
